$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B, shifting the existing "Valor" column to C.
$ws.Columns("B:B").Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("D1").Value = "Colocação"

# Common variable label used on every data row
$varLabel = "Trabalho como origem na renda (%): média de 2012 a 2019"

$ws.Range("B2").Value = $varLabel
$ws.Range("B3").Value = $varLabel
$ws.Range("B4").Value = $varLabel
$ws.Range("B5").Value = $varLabel
$ws.Range("B6").Value = $varLabel
$ws.Range("B7").Value = $varLabel
$ws.Range("B8").Value = $varLabel
$ws.Range("B9").Value = $varLabel
$ws.Range("B10").Value = $varLabel

# Ranking column (only populated for the first eight data rows)
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "23º"
